# Apply the "data/types.xlsx" edit:
#  - B9: 0 -> 1
#  - Clear out the old "nano Ti" / "VF4-Ti3C2" / "MBH-*" / "NdF3" / "TiF3" / "V2C"
#    rows (A34:J48) - these rows' content is removed (row 34 and rows 36-48
#    disappear entirely, row 35 keeps only its pre-existing formatting with
#    no value), which also prunes the now-unused shared strings.
#  - Update the current selection to N28, matching the saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 1

$ws.Range("A34:J48").ClearContents()

$ws.Range("N28").Select()
